# Insert three new columns (Model_Family, Engine, Split_Type) after "Model"
# and before "Source", then re-populate the Source/metric columns and
# regroup the data rows to reflect the new Model_Family / Engine / Split_Type
# breakdown (mirrors the OOXML diff applied to sheet1.xml / sharedStrings.xml).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 columns before column C (old "Avg_AIC"), which pushes the old
# "Source" header (col B) to stay at B, and shifts everything from C onward
# three columns to the right. We actually want the new columns right after
# "Model" (col A) and before "Source" (old col B), so insert at column B.
$ws.Range("B:D").Insert() | Out-Null

# Header row
$ws.Range("A1").Value = "Model"
$ws.Range("B1").Value = "Model_Family"
$ws.Range("C1").Value = "Engine"
$ws.Range("D1").Value = "Split_Type"
$ws.Range("E1").Value = "Source"
$ws.Range("F1").Value = "Avg_AIC"
$ws.Range("G1").Value = "Avg_BIC"
$ws.Range("H1").Value = "Avg_LogLik"
$ws.Range("I1").Value = "Avg_MSE"
$ws.Range("J1").Value = "Avg_MAE"

# Data rows: Model, Model_Family, Engine, Split_Type, Source, Avg_AIC, Avg_BIC, Avg_LogLik, Avg_MSE, Avg_MAE
$data = @(
    @("TGARCH",      "GARCH",    "N/A", "TS_CV", "Time_Series_CV", -6.23674331462475,  -6.17773880124684,  1566.18582865619,  0.000265231788437081, 0.00997959312049119),
    @("eGARCH",       "GARCH",    "N/A", "TS_CV", "Time_Series_CV", -6.16499588686336,  -6.10599137348545,  1548.24897171584,  0.000265255247246037, 0.00998057275133209),
    @("eGARCH",       "NF-GARCH", "N/A", "N/A",   "NF-GARCH",        27401.3105107422,   27439.8014735373,  -13694.6552553711, 0,                     0),
    @("fGARCH",       "NF-GARCH", "N/A", "N/A",   "NF-GARCH",       -28273.7513163331,  -28235.260353538,   14142.8756581666,  0,                     0),
    @("gjrGARCH",     "GARCH",    "N/A", "TS_CV", "Time_Series_CV", -6.23529221703763,  -6.17628770365972,  1565.82305425941,  0.000264427773901888, 0.0099650836905736),
    @("gjrGARCH",     "NF-GARCH", "N/A", "N/A",   "NF-GARCH",       -28244.7567895128,  -28206.2658267177,  14128.3783947564,  0,                     0),
    @("sGARCH",       "NF-GARCH", "N/A", "N/A",   "NF-GARCH",       -27993.3295905363,  -27964.46136844,    14001.1647952681,  0,                     0),
    @("sGARCH_norm",  "GARCH",    "N/A", "TS_CV", "Time_Series_CV", -6.14008516085673,  -6.10636829606935,  1539.02129021418,  0.00026524306076551,  0.00997754931432829),
    @("sGARCH_sstd",  "GARCH",    "N/A", "TS_CV", "Time_Series_CV", -6.22901892609995,  -6.17844362891889,  1563.25473152499,  0.000265218411076707, 0.00997773042969076)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $ws.Cells.Item($row, 9).Value = $r[8]
    $ws.Cells.Item($row, 10).Value = $r[9]
    $row++
}
